# SANTA_CRUZ_DO_SUL.xlsx - automatic update
# 1) Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Remove the "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
